$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 2")

# --- Row 13: fill in the new "EF Core valmis" time-log entry ---

# B13 is formatted as a date-like text column (dd/mm/yyyy;@) but the existing
# entries really hold plain text strings (e.g. "10.02.2020") rather than real
# dates, so force text formatting while assigning to avoid Excel silently
# converting the literal into a date serial number.
$b13Format = $ws.Range("B13").NumberFormat()
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "10.02.2020"
$ws.Range("B13").NumberFormat = $b13Format

$ws.Range("C13").Value = 0.77083333333333337
$ws.Range("D13").Value = 0.92013888888888884
$ws.Range("F13").Value = 215
$ws.Range("G13").Value = "Koduneül"
$ws.Range("H13").Value = "EF Core valmis"

# --- Row 14: the "podcast" activity note is cleared out ---
$ws.Range("G14").ClearContents()

# Move the active selection to G14, matching the saved view state
$ws.Range("G14").Select()

$wb.Save()
